$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header row: A1 was "Gene" -> "Target", B1 was "Primer" -> "Gene", C1 stays "HKG"
$ws.Range("A1").Value = "Target"
$ws.Range("B1").Value = "Gene"
$ws.Range("C1").Value = "HKG"

# Convert the HKG boolean flags (TRUE/FALSE) into text flags ("Y"/"N")
$ws.Range("C2:C4").Value = "Y"
$ws.Range("C5:C17").Value = "N"

# Move the active selection to F7, matching the saved view state
$ws.Range("F7").Select()
